$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 37.99343233333334
$ws.Range("H2").Value = 113.980297
$ws.Range("I2").Value = 0.3685480664467733
$ws.Range("J2").Value = 0.3685480664467734
$ws.Range("M2").Value = 91.74689966666665
$ws.Range("N2").Value = 275.2406989999999
$ws.Range("O2").Value = 0.1908387282982634
$ws.Range("P2").Value = 0.1908387282982634
$ws.Range("Q2").Value = 3485.779624278623
$ws.Range("R2").Value = 31372.0166185076
$ws.Range("S2").Value = 0.0703332443174861
$ws.Range("T2").Value = 0.07033324431748611

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 37.99343233333334
$ws.Range("H3").Value = 113.980297
$ws.Range("I3").Value = 0.3685480664467733
$ws.Range("J3").Value = 0.3685480664467734
$ws.Range("O3").Value = 0.296899627499751
$ws.Range("P3").Value = 0.296899627499751
$ws.Range("Q3").Value = 5423.043222008112
$ws.Range("R3").Value = 48807.388998073
$ws.Range("S3").Value = 0.1094217836438005
$ws.Range("T3").Value = 0.1094217836438005

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 37.99343233333334
$ws.Range("H4").Value = 113.980297
$ws.Range("I4").Value = 0.3685480664467733
$ws.Range("J4").Value = 0.3685480664467734
$ws.Range("M4").Value = 167.6324513333334
$ws.Range("N4").Value = 502.8973540000001
$ws.Range("O4").Value = 0.348684957750095
$ws.Range("P4").Value = 0.348684957750095
$ws.Range("Q4").Value = 6368.932196603795
$ws.Range("R4").Value = 57320.38976943415
$ws.Range("S4").Value = 0.1285071669778724
$ws.Range("T4").Value = 0.1285071669778724

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 37.99343233333334
$ws.Range("H5").Value = 113.980297
$ws.Range("I5").Value = 0.3685480664467733
$ws.Range("J5").Value = 0.3685480664467734
$ws.Range("M5").Value = 78.64050433333334
$ws.Range("N5").Value = 235.921513
$ws.Range("O5").Value = 0.1635766864518907
$ws.Range("P5").Value = 0.1635766864518907
$ws.Range("Q5").Value = 2987.822680047707
$ws.Range("R5").Value = 26890.40412042936
$ws.Range("S5").Value = 0.06028587150761442
$ws.Range("T5").Value = 0.06028587150761443

# Row 6
$ws.Range("I6").Value = 0.3751865155371754
$ws.Range("J6").Value = 0.3751865155371755
$ws.Range("M6").Value = 91.74689966666665
$ws.Range("N6").Value = 275.2406989999999
$ws.Range("O6").Value = 0.1908387282982634
$ws.Range("P6").Value = 0.1908387282982634
$ws.Range("Q6").Value = 3548.567012635403
$ws.Range("R6").Value = 31937.10311371863
$ws.Range("S6").Value = 0.07160011749977117
$ws.Range("T6").Value = 0.07160011749977122

# Row 7
$ws.Range("I7").Value = 0.3751865155371754
$ws.Range("J7").Value = 0.3751865155371755
$ws.Range("O7").Value = 0.296899627499751
$ws.Range("P7").Value = 0.296899627499751
$ws.Range("S7").Value = 0.1113927367059169
$ws.Range("T7").Value = 0.1113927367059169

# Row 8
$ws.Range("I8").Value = 0.3751865155371754
$ws.Range("J8").Value = 0.3751865155371755
$ws.Range("M8").Value = 167.6324513333334
$ws.Range("N8").Value = 502.8973540000001
$ws.Range("O8").Value = 0.348684957750095
$ws.Range("P8").Value = 0.348684957750095
$ws.Range("Q8").Value = 6483.65219108105
$ws.Range("R8").Value = 58352.86971972944
$ws.Range("S8").Value = 0.1308218943184853
$ws.Range("T8").Value = 0.1308218943184854

# Row 9
$ws.Range("I9").Value = 0.3751865155371754
$ws.Range("J9").Value = 0.3751865155371755
$ws.Range("M9").Value = 78.64050433333334
$ws.Range("N9").Value = 235.921513
$ws.Range("O9").Value = 0.1635766864518907
$ws.Range("P9").Value = 0.1635766864518907
$ws.Range("Q9").Value = 3041.640649963742
$ws.Range("R9").Value = 27374.76584967368
$ws.Range("S9").Value = 0.06137176701300195
$ws.Range("T9").Value = 0.06137176701300197

# Row 10
$ws.Range("G10").Value = 26.18781466666667
$ws.Range("H10").Value = 78.563444
$ws.Range("I10").Value = 0.2540299169390597
$ws.Range("J10").Value = 0.2540299169390597
$ws.Range("M10").Value = 91.74689966666665
$ws.Range("N10").Value = 275.2406989999999
$ws.Range("O10").Value = 0.1908387282982634
$ws.Range("P10").Value = 0.1908387282982634
$ws.Range("Q10").Value = 2402.650804711928
$ws.Range("R10").Value = 21623.85724240735
$ws.Range("S10").Value = 0.04847874629836362
$ws.Range("T10").Value = 0.04847874629836364

# Row 11
$ws.Range("G11").Value = 26.18781466666667
$ws.Range("H11").Value = 78.563444
$ws.Range("I11").Value = 0.2540299169390597
$ws.Range("J11").Value = 0.2540299169390597
$ws.Range("O11").Value = 0.296899627499751
$ws.Range("P11").Value = 0.296899627499751
$ws.Range("Q11").Value = 3737.952643532889
$ws.Range("R11").Value = 33641.573791796
$ws.Range("S11").Value = 0.0754213877129995
$ws.Range("T11").Value = 0.07542138771299951

# Row 12
$ws.Range("G12").Value = 26.18781466666667
$ws.Range("H12").Value = 78.563444
$ws.Range("I12").Value = 0.2540299169390597
$ws.Range("J12").Value = 0.2540299169390597
$ws.Range("M12").Value = 167.6324513333334
$ws.Range("N12").Value = 502.8973540000001
$ws.Range("O12").Value = 0.348684957750095
$ws.Range("P12").Value = 0.348684957750095
$ws.Range("Q12").Value = 4389.927567636354
$ws.Range("R12").Value = 39509.34810872719
$ws.Range("S12").Value = 0.08857641085515618
$ws.Range("T12").Value = 0.08857641085515619

# Row 13
$ws.Range("G13").Value = 26.18781466666667
$ws.Range("H13").Value = 78.563444
$ws.Range("I13").Value = 0.2540299169390597
$ws.Range("J13").Value = 0.2540299169390597
$ws.Range("M13").Value = 78.64050433333334
$ws.Range("N13").Value = 235.921513
$ws.Range("O13").Value = 0.1635766864518907
$ws.Range("P13").Value = 0.1635766864518907
$ws.Range("Q13").Value = 2059.422952774531
$ws.Range("R13").Value = 18534.80657497077
$ws.Range("S13").Value = 0.0415533720725404
$ws.Range("T13").Value = 0.04155337207254041

# Row 14
$ws.Range("G14").Value = 0.2304566666666667
$ws.Range("H14").Value = 0.69137
$ws.Range("I14").Value = 0.002235501076991453
$ws.Range("J14").Value = 0.002235501076991454
$ws.Range("M14").Value = 91.74689966666665
$ws.Range("N14").Value = 275.2406989999999
$ws.Range("O14").Value = 0.1908387282982634
$ws.Range("P14").Value = 0.1908387282982634
$ws.Range("Q14").Value = 21.14368467418111
$ws.Range("R14").Value = 190.29316206763
$ws.Range("S14").Value = 0.0004266201826424471
$ws.Range("T14").Value = 0.0004266201826424472

# Row 15
$ws.Range("G15").Value = 0.2304566666666667
$ws.Range("H15").Value = 0.69137
$ws.Range("I15").Value = 0.002235501076991453
$ws.Range("J15").Value = 0.002235501076991454
$ws.Range("O15").Value = 0.296899627499751
$ws.Range("P15").Value = 0.296899627499751
$ws.Range("Q15").Value = 32.89453959222222
$ws.Range("R15").Value = 296.05085633
$ws.Range("S15").Value = 0.0006637194370340546
$ws.Range("T15").Value = 0.0006637194370340547

# Row 16
$ws.Range("G16").Value = 0.2304566666666667
$ws.Range("H16").Value = 0.69137
$ws.Range("I16").Value = 0.002235501076991453
$ws.Range("J16").Value = 0.002235501076991454
$ws.Range("M16").Value = 167.6324513333334
$ws.Range("N16").Value = 502.8973540000001
$ws.Range("O16").Value = 0.348684957750095
$ws.Range("P16").Value = 0.348684957750095
$ws.Range("Q16").Value = 38.63201595944223
$ws.Range("R16").Value = 347.6881436349801
$ws.Range("S16").Value = 0.0007794855985810567
$ws.Range("T16").Value = 0.0007794855985810569

# Row 17
$ws.Range("G17").Value = 0.2304566666666667
$ws.Range("H17").Value = 0.69137
$ws.Range("I17").Value = 0.002235501076991453
$ws.Range("J17").Value = 0.002235501076991454
$ws.Range("M17").Value = 78.64050433333334
$ws.Range("N17").Value = 235.921513
$ws.Range("O17").Value = 0.1635766864518907
$ws.Range("P17").Value = 0.1635766864518907
$ws.Range("Q17").Value = 18.12322849364556
$ws.Range("R17").Value = 163.10905644281
$ws.Range("S17").Value = 0.0003656758587338949
$ws.Range("T17").Value = 0.000365675858733895
